# Update the "Förändrad" (Changed) date column (column C) for every data
# row on the active sheet: bump the stored date serial number by one day
# (2023-10-03 -> 2023-10-04), leaving every other cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -ne $null) {
        $cell.Value = $current + 1
    }
}
